# Add rural communities interactions back in, fix de_dg files
# The underlying edit: row 12 ("Legislature") is removed from Sheet1,
# which shifts all subsequent rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire row 12 (the "Legislature" row), shifting rows 13:29 up to 12:28
$ws.Rows.Item(12).Delete()

# Update the view/selection state to reflect the post-delete state seen in the file
$ws.Activate()
$ws.Range("A12:XFD12").Select()
